$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1.35427148362301
$ws.Range("C2").Value = 0.5399220165771437
$ws.Range("D2").Value = 0
